$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> (new DAMSLTag, new DialogAct)
$changes = @{
    7   = @("aa", "Agree/Accept")
    8   = @("sd", "Statement-non-opinion")
    19  = @("b", "Acknowledge (Backchannel)")
    35  = @("sv", "Statement-opinion")
    45  = @("b", "Acknowledge (Backchannel)")
    49  = @("aa", "Agree/Accept")
    55  = @("sv", "Statement-opinion")
    65  = @("b", "Acknowledge (Backchannel)")
    127 = @("aa", "Agree/Accept")
    134 = @("sd", "Statement-non-opinion")
    138 = @("sd", "Statement-non-opinion")
    140 = @("%", "Uninterpretable")
    175 = @("sv", "Statement-opinion")
    176 = @("aa", "Agree/Accept")
    179 = @("b", "Acknowledge (Backchannel)")
    187 = @("sd", "Statement-non-opinion")
    193 = @("ba", "Appreciation")
    197 = @("ba", "Appreciation")
    204 = @("b", "Acknowledge (Backchannel)")
    225 = @("b", "Acknowledge (Backchannel)")
    230 = @("aa", "Agree/Accept")
    235 = @("b", "Acknowledge (Backchannel)")
    238 = @("aa", "Agree/Accept")
    248 = @("sd", "Statement-non-opinion")
    258 = @("sd", "Statement-non-opinion")
    262 = @("sv", "Statement-opinion")
    278 = @("sd", "Statement-non-opinion")
    279 = @("%", "Uninterpretable")
    281 = @("b", "Acknowledge (Backchannel)")
    284 = @("sd", "Statement-non-opinion")
    290 = @("sd", "Statement-non-opinion")
}

foreach ($row in $changes.Keys) {
    $values = $changes[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
